# Update E8 from "Good Morning" to "GIT UPDATE", matching the commit's
# "update file with jgit" change. This removes the old shared string
# ("Good Morning") and appends a new one ("GIT UPDATE") at the end of
# the shared-strings table, which is exactly what happens when the
# only remaining reference to a string is overwritten with new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

# Leave the selection on the edited cell, as recorded in the saved file.
$ws.Range("E8").Select()
